{"js": "// Paraphrase two short spans inside the \"Second, for every single stock...\" bullet:\n//   1) \"every single\" -> \"the user-selected\"\n//   2) the final \"...the model attempts to predict a stock's closing price.\"\n//      -> \"...the model attempts to predict the stock's closing price.\"\n// (there is an earlier, untouched \"attempts to predict a stock's closing price\"\n//  in the same paragraph, so we must only touch the last occurrence).\n\n// --- Change 1: \"every single\" -> \"the user-selected\" -------------------\nlet results = context.document.body.search(\"every single\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"the user-selected\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: last \"predict a stock's closing price.\" -> \"...the stock's...\" ----\n// Use enough surrounding context to uniquely identify the final occurrence\n// (the paragraph also contains an earlier, unrelated \"predict a stock's\n// closing price\" that must stay exactly as-is).\nconst target = \"the model attempts to predict a stock's closing price. \";\nconst replacement = \"the model attempts to predict the stock's closing price. \";\n\nlet results2 = context.document.body.search(target, { matchCase: true, matchWholeWord: false });\nresults2.load(\"items\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  results2.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Paraphrase two short spans inside the \"Second, for every single stock...\"\n# bullet point:\n#   1) \"every single\"  -> \"the user-selected\"\n#   2) the final occurrence of \"...the model attempts to predict a stock's\n#      closing price.\" -> \"...the model attempts to predict the stock's\n#      closing price.\"\n# (an earlier, unrelated \"attempts to predict a stock's closing price\" in the\n#  very same paragraph must stay untouched, so change 2 uses enough\n#  surrounding context to uniquely identify only the last occurrence.)\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"every single\" -> \"the user-selected\" ----------------------\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.ClearFormatting()\n$find1.Text = \"every single\"\n$find1.MatchCase = $true\n$found1 = $find1.Execute()\nif ($found1) {\n    $range1.Text = \"the user-selected\"\n}\n\n# --- Change 2: last \"...predict a stock's closing price.\" -> \"...the stock's...\" ----\n$target2 = \"the model attempts to predict a stock's closing price. \"\n$replacement2 = \"the model attempts to predict the stock's closing price. \"\n\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = $target2\n$find2.MatchCase = $true\n$found2 = $find2.Execute()\nif ($found2) {\n    $range2.Text = $replacement2\n}\n"}
